$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.674.55'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.959.67'
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '244.66'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  +2.33%  '
$ws.Range("D7").Value = '61.79'
$ws.Range("E7").Value = '  +7.94%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +4.75%  '
$ws.Range("D10").Value = '0.0794'
$ws.Range("E10").Value = '  -6.48%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").Value = '14.27'
$ws.Range("E12").Value = '  +5.96%  '
$ws.Range("D13").Value = '22.08'
$ws.Range("E13").Value = '  +3.66%  '
$ws.Range("D14").Value = '0.833'
$ws.Range("E14").Value = '  +2.82%  '
$ws.Range("D15").Value = '2.245.10'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '5.30'
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").Value = '1.959.61'
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("D18").Value = '36.548.77'
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '69.74'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").Value = '0.0₃0853'
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = '230.12'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("D22").Value = '5.08'
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  +5.19%  '
$ws.Range("E25").Value = '  +2.89%  '
$ws.Range("E26").Value = '  +7.50%  '
$ws.Range("D27").Value = '9.19'
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").Value = '160.53'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = '19.42'
$ws.Range("E29").Value = '  +1.07%  '
$ws.Range("E30").Value = '  +18.10%  '
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("D32").Value = '4.77'
$ws.Range("E32").Value = '  +4.65%  '
$ws.Range("D33").Value = '0.0616'
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = '4.46'
$ws.Range("E34").Value = '  +7.04%  '
$ws.Range("D35").Value = '3.56'
$ws.Range("E35").Value = '  +13.96%  '
$ws.Range("D36").Value = '2.28'
$ws.Range("E36").Value = '  +4.88%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").Value = '5.52'
$ws.Range("E39").Value = '  -9.40%  '
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").Value = '1.18'
$ws.Range("E42").Value = '  +2.21%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").Value = '16.03'
$ws.Range("E44").Value = '  +2.40%  '
$ws.Range("D45").Value = '1.369.33'
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("D46").Value = '88.66'
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("D48").Value = '7.14'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("D50").Value = '45.26'
$ws.Range("E50").Value = '  +5.14%  '
$ws.Range("D51").Value = '2.125.93'
$ws.Range("E51").Value = '  +0.47%  '
